$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 1 (BTec_Logo-Orange): image1.jpg -> image2.jpg
$sec.Headers.Item(2).Shapes.Item(1).Name = "image2.jpg"

# Header 2 (BTec_Logo-Orange): image1.jpg -> image2.jpg
$sec.Headers.Item(1).Shapes.Item(1).Name = "image2.jpg"

# Footer 1 (PearsonLogo): image2.png -> image1.png
$sec.Footers.Item(2).Shapes.Item(1).Name = "image1.png"

# Footer 2 (PearsonLogo): image2.png -> image1.png
$sec.Footers.Item(1).Shapes.Item(1).Name = "image1.png"
